$d = $word.ActiveDocument

# 1. Remove the title paragraph ("Uso de Servicios de AWS en la Migración
#    de una Panadería a la Nube") entirely, including its paragraph mark.
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Delete()

# 2. In the (now first) subtitle paragraph, drop the "Tabla 1" run and the
#    line break before "Aplicación de los Servicios de AWS en una Panadería".
$d.Content.Find.Execute("Tabla 1^l", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2)

# 3. Merge the "Uso " / "de servicio" runs in the header row of the table
#    into a single run "Uso de servicio" (keeping the bold formatting that
#    both original runs already shared).
$table = $d.Tables.Item(1)
$cell = $table.Cell(1, 2)
$cell.Range.Find.Execute("Uso de servicio", $true, $false, $false, $false, `
                          $false, $true, 1, $false, "Uso de servicio", 2)

# 4. Register the two new paragraph styles used for table captions/content.
$contenido = $d.Styles.Add("Contenidodelatabla", 1)
$contenido.NameLocal = "Contenido de la tabla"
$contenido.BaseStyle = "Normal"
$contenido.QuickStyle = $true
$contenido.ParagraphFormat.WidowControl = $false
$contenido.ParagraphFormat.NoLineNumber = $true

$titulo = $d.Styles.Add("Ttulodelatabla", 1)
$titulo.NameLocal = "Título de la tabla"
$titulo.BaseStyle = "Contenidodelatabla"
$titulo.QuickStyle = $true
$titulo.ParagraphFormat.NoLineNumber = $true
$titulo.ParagraphFormat.Alignment = 1
$titulo.Font.Bold = $true
$titulo.Font.BoldBi = $true
